# Update project dates to 2025-08-01 .. 2025-11-20(ish) and format dates
# without a time component (new date-only number format, numFmtId 165).
#
# Sheet "Sprint Backlog 1": column G (rows 4-93) holds the "ESTIMATED DATE
# OF TASK COMPLETION" for each task. The 90 tasks are regrouped into
# batches of 3 tasks/day starting 2025-08-01 (serial 45870).
#
# Sheet "SB BD Ch01": burndown chart data (rows 3-40) is regenerated to
# span the new 38-day project window (2025-08-01 .. 2025-09-07), linearly
# decreasing the task count from 90 to 0 and rounding to 1 decimal place.

$wb = $excel.ActiveWorkbook

$wsBacklog = $wb.Worksheets.Item("Sprint Backlog 1")
$wsBurndown = $wb.Worksheets.Item("SB BD Ch01")

# --- Sprint Backlog 1: column G (rows 4..93) ---------------------------
$backlogFirstRow = 4
$backlogTaskCount = 90
$backlogStartDate = 45870
$tasksPerDay = 3

$rangeG = $wsBacklog.Range("G4:G93")
$rangeG.NumberFormat = "yyyy-mm-dd"

for ($i = 0; $i -lt $backlogTaskCount; $i++) {
    $row = $backlogFirstRow + $i
    $dayOffset = [math]::Floor($i / $tasksPerDay)
    $date = $backlogStartDate + $dayOffset
    $wsBacklog.Cells.Item($row, 7).Value = $date
}

# --- SB BD Ch01: burndown chart rows 3..40 -----------------------------
$burndownFirstRow = 3
$burndownDayCount = 38
$burndownStartDate = 45870
$burndownStartTasks = 90

$rangeA = $wsBurndown.Range("A3:A40")
$rangeA.NumberFormat = "yyyy-mm-dd"

for ($i = 0; $i -lt $burndownDayCount; $i++) {
    $row = $burndownFirstRow + $i
    $date = $burndownStartDate + $i
    $remaining = [math]::Round($burndownStartTasks - ($i * $burndownStartTasks / ($burndownDayCount - 1)), 1)

    $wsBurndown.Cells.Item($row, 1).Value = $date
    $wsBurndown.Cells.Item($row, 2).Value = $remaining
    $wsBurndown.Cells.Item($row, 3).Value = $remaining
}
